$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows of data (2-9, excluding the header row 1) got re-ordered/shuffled
# (weekly data re-sequenced). Row 8 is unchanged. Apply the new values for
# columns D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), Q (Unidad de comercializacion) and
# S (Precio $/Kg) for each affected row.

# Row 2
$ws.Range("D2").Value = Get-Date -Year 2021 -Month 4 -Day 23 -Hour 0 -Minute 0 -Second 0
$ws.Range("M2").Value = 300
$ws.Range("N2").Value = 7000
$ws.Range("O2").Value = 7000
$ws.Range("P2").Value = 7000
$ws.Range("Q2").Value = '$/caja 14 kilos empedrada'
$ws.Range("S2").Value = 500

# Row 3
$ws.Range("D3").Value = Get-Date -Year 2021 -Month 7 -Day 23 -Hour 0 -Minute 0 -Second 0
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 10000
$ws.Range("O3").Value = 10000
$ws.Range("P3").Value = 10000
$ws.Range("Q3").Value = '$/caja 14 kilos'
$ws.Range("S3").Value = 714

# Row 4
$ws.Range("D4").Value = Get-Date -Year 2020 -Month 12 -Day 11 -Hour 0 -Minute 0 -Second 0
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 7000
$ws.Range("O4").Value = 7000
$ws.Range("P4").Value = 7000
$ws.Range("Q4").Value = '$/caja 14 kilos empedrada'
$ws.Range("S4").Value = 500

# Row 5
$ws.Range("D5").Value = Get-Date -Year 2021 -Month 7 -Day 20 -Hour 0 -Minute 0 -Second 0
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 11000
$ws.Range("O5").Value = 11000
$ws.Range("P5").Value = 11000
$ws.Range("Q5").Value = '$/caja 14 kilos'
$ws.Range("S5").Value = 786

# Row 6
$ws.Range("D6").Value = Get-Date -Year 2021 -Month 10 -Day 22 -Hour 0 -Minute 0 -Second 0
$ws.Range("M6").Value = 180
$ws.Range("N6").Value = 9000
$ws.Range("O6").Value = 9000
$ws.Range("P6").Value = 9000
$ws.Range("Q6").Value = '$/caja 14 kilos empedrada'
$ws.Range("S6").Value = 643

# Row 7
$ws.Range("D7").Value = Get-Date -Year 2020 -Month 11 -Day 27 -Hour 0 -Minute 0 -Second 0
$ws.Range("M7").Value = 120
$ws.Range("N7").Value = 7000
$ws.Range("O7").Value = 7000
$ws.Range("P7").Value = 7000
$ws.Range("Q7").Value = '$/caja 14 kilos empedrada'
$ws.Range("S7").Value = 500

# Row 8 unchanged

# Row 9
$ws.Range("D9").Value = Get-Date -Year 2021 -Month 6 -Day 4 -Hour 0 -Minute 0 -Second 0
$ws.Range("M9").Value = 300
$ws.Range("N9").Value = 10000
$ws.Range("O9").Value = 10000
$ws.Range("P9").Value = 10000
$ws.Range("Q9").Value = '$/caja 14 kilos empedrada'
$ws.Range("S9").Value = 714
